$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column C ("Station") entries to stay text ("48") instead of being
# auto-coerced to a number, matching the source data's inline-string type.
$ws.Range("C5:C6").NumberFormat = "@"

# Row 5: SOLEMON2024 | ITA17 | 48 | 2-RAP | RAJAAST | 1 | -1 | SIMRANDO | Y
$ws.Range("A5").Value = "SOLEMON2024"
$ws.Range("B5").Value = "ITA17"
$ws.Range("C5").Value = "48"
$ws.Range("D5").Value = "2-RAP"
$ws.Range("E5").Value = "RAJAAST"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = -1
$ws.Range("H5").Value = "SIMRANDO"
$ws.Range("I5").Value = "Y"

# Row 6: SOLEMON2024 | ITA17 | 48 | 1-RAP | RAJAAST | 1 | -1 | SIMRANDO | Y
$ws.Range("A6").Value = "SOLEMON2024"
$ws.Range("B6").Value = "ITA17"
$ws.Range("C6").Value = "48"
$ws.Range("D6").Value = "1-RAP"
$ws.Range("E6").Value = "RAJAAST"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = -1
$ws.Range("H6").Value = "SIMRANDO"
$ws.Range("I6").Value = "Y"

# Restore default (General/no-style) formatting now that the values are
# locked in as text, so the new rows don't carry a lingering style index.
$ws.Range("C5:C6").Style = "Normal"
